$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Arf1"
$ws.Range("C2").Value = "Chrm3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 52.138213
$ws.Range("H2").Value = 156.414639
$ws.Range("I2").Value = 0.2220849502516424
$ws.Range("J2").Value = 0.2220849502516423
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4727823333333334
$ws.Range("N2").Value = 1.418347
$ws.Range("O2").Value = 0.4785956754713925
$ws.Range("P2").Value = 0.4785956754713925
$ws.Range("Q2").Value = 24.65002599797034
$ws.Range("R2").Value = 221.850233981733
$ws.Range("S2").Value = 0.1062888967777154
$ws.Range("T2").Value = 0.1062888967777154

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Arf1"
$ws.Range("C3").Value = "Chrm3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 52.138213
$ws.Range("H3").Value = 156.414639
$ws.Range("I3").Value = 0.2220849502516424
$ws.Range("J3").Value = 0.2220849502516423
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1088336666666667
$ws.Range("N3").Value = 0.326501
$ws.Range("O3").Value = 0.1101718878645953
$ws.Range("P3").Value = 0.1101718878645953
$ws.Range("Q3").Value = 5.674392894237666
$ws.Range("R3").Value = 51.06953604813899
$ws.Range("S3").Value = 0.02446751823553817
$ws.Range("T3").Value = 0.02446751823553816

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Arf1"
$ws.Range("C4").Value = "Chrm3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 52.138213
$ws.Range("H4").Value = 156.414639
$ws.Range("I4").Value = 0.2220849502516424
$ws.Range("J4").Value = 0.2220849502516423
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4062373333333333
$ws.Range("N4").Value = 1.218712
$ws.Range("O4").Value = 0.4112324366640122
$ws.Range("P4").Value = 0.4112324366640122
$ws.Range("Q4").Value = 21.18048861388533
$ws.Range("R4").Value = 190.624397524968
$ws.Range("S4").Value = 0.09132853523838882
$ws.Range("T4").Value = 0.0913285352383888

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Arf1"
$ws.Range("C5").Value = "Chrm3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 67.324
$ws.Range("H5").Value = 201.972
$ws.Range("I5").Value = 0.2867694600645705
$ws.Range("J5").Value = 0.2867694600645705
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4727823333333334
$ws.Range("N5").Value = 1.418347
$ws.Range("O5").Value = 0.4785956754713925
$ws.Range("P5").Value = 0.4785956754713925
$ws.Range("Q5").Value = 31.82959780933333
$ws.Range("R5").Value = 286.466380284
$ws.Range("S5").Value = 0.1372466234441696
$ws.Range("T5").Value = 0.1372466234441696

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Arf1"
$ws.Range("C6").Value = "Chrm3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 67.324
$ws.Range("H6").Value = 201.972
$ws.Range("I6").Value = 0.2867694600645705
$ws.Range("J6").Value = 0.2867694600645705
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1088336666666667
$ws.Range("N6").Value = 0.326501
$ws.Range("O6").Value = 0.1101718878645953
$ws.Range("P6").Value = 0.1101718878645953
$ws.Range("Q6").Value = 7.327117774666666
$ws.Range("R6").Value = 65.944059972
$ws.Range("S6").Value = 0.0315939327972244
$ws.Range("T6").Value = 0.03159393279722439

# Row 7: FAPs -> MuSCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Arf1"
$ws.Range("C7").Value = "Chrm3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 67.324
$ws.Range("H7").Value = 201.972
$ws.Range("I7").Value = 0.2867694600645705
$ws.Range("J7").Value = 0.2867694600645705
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4062373333333333
$ws.Range("N7").Value = 1.218712
$ws.Range("O7").Value = 0.4112324366640122
$ws.Range("P7").Value = 0.4112324366640122
$ws.Range("Q7").Value = 27.34952222933333
$ws.Range("R7").Value = 246.145700064
$ws.Range("S7").Value = 0.1179289038231765
$ws.Range("T7").Value = 0.1179289038231764

# Row 8: MuSCs -> ECs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Arf1"
$ws.Range("C8").Value = "Chrm3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 64.999789
$ws.Range("H8").Value = 194.999367
$ws.Range("I8").Value = 0.2768693838132169
$ws.Range("J8").Value = 0.2768693838132169
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.4727823333333334
$ws.Range("N8").Value = 1.418347
$ws.Range("O8").Value = 0.4785956754713925
$ws.Range("P8").Value = 0.4785956754713925
$ws.Range("Q8").Value = 30.73075190959434
$ws.Range("R8").Value = 276.576767186349
$ws.Range("S8").Value = 0.1325084897634348
$ws.Range("T8").Value = 0.1325084897634347

# Row 9: MuSCs -> FAPs
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Arf1"
$ws.Range("C9").Value = "Chrm3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 64.999789
$ws.Range("H9").Value = 194.999367
$ws.Range("I9").Value = 0.2768693838132169
$ws.Range("J9").Value = 0.2768693838132169
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1088336666666667
$ws.Range("N9").Value = 0.326501
$ws.Range("O9").Value = 0.1101718878645953
$ws.Range("P9").Value = 0.1101718878645953
$ws.Range("Q9").Value = 7.074165369429667
$ws.Range("R9").Value = 63.667488324867
$ws.Range("S9").Value = 0.03050322270660932
$ws.Range("T9").Value = 0.03050322270660932

# Row 10: MuSCs -> MuSCs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Arf1"
$ws.Range("C10").Value = "Chrm3"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 64.999789
$ws.Range("H10").Value = 194.999367
$ws.Range("I10").Value = 0.2768693838132169
$ws.Range("J10").Value = 0.2768693838132169
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4062373333333333
$ws.Range("N10").Value = 1.218712
$ws.Range("O10").Value = 0.4112324366640122
$ws.Range("P10").Value = 0.4112324366640122
$ws.Range("Q10").Value = 26.40534095058934
$ws.Range("R10").Value = 237.648068555304
$ws.Range("S10").Value = 0.1138576713431728
$ws.Range("T10").Value = 0.1138576713431728

# Row 11: Resolving-Mac -> ECs
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Arf1"
$ws.Range("C11").Value = "Chrm3"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 50.30497766666667
$ws.Range("H11").Value = 150.914933
$ws.Range("I11").Value = 0.2142762058705703
$ws.Range("J11").Value = 0.2142762058705703
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.4727823333333334
$ws.Range("N11").Value = 1.418347
$ws.Range("O11").Value = 0.4785956754713925
$ws.Range("P11").Value = 0.4785956754713925
$ws.Range("Q11").Value = 23.78330471952789
$ws.Range("R11").Value = 214.049742475751
$ws.Range("S11").Value = 0.1025516654860728
$ws.Range("T11").Value = 0.1025516654860728

# Row 12: Resolving-Mac -> FAPs
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Arf1"
$ws.Range("C12").Value = "Chrm3"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 50.30497766666667
$ws.Range("H12").Value = 150.914933
$ws.Range("I12").Value = 0.2142762058705703
$ws.Range("J12").Value = 0.2142762058705703
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.1088336666666667
$ws.Range("N12").Value = 0.326501
$ws.Range("O12").Value = 0.1101718878645953
$ws.Range("P12").Value = 0.1101718878645953
$ws.Range("Q12").Value = 5.474875171048111
$ws.Range("R12").Value = 49.273876539433
$ws.Range("S12").Value = 0.02360721412522341
$ws.Range("T12").Value = 0.02360721412522341

# Row 13: Resolving-Mac -> MuSCs
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Arf1"
$ws.Range("C13").Value = "Chrm3"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 50.30497766666667
$ws.Range("H13").Value = 150.914933
$ws.Range("I13").Value = 0.2142762058705703
$ws.Range("J13").Value = 0.2142762058705703
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.4062373333333333
$ws.Range("N13").Value = 1.218712
$ws.Range("O13").Value = 0.4112324366640122
$ws.Range("P13").Value = 0.4112324366640122
$ws.Range("Q13").Value = 20.43575998069955
$ws.Range("R13").Value = 183.921839826296
$ws.Range("S13").Value = 0.08811732625927415
$ws.Range("T13").Value = 0.08811732625927414

$ws.Rows("14:17").Delete()